$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 27739.75
$ws.Range("J57").Value = 27739.75
$ws.Range("L57").Value = 83219.25
$ws.Range("N57").Value = -84217.25
$ws.Range("H76").Value = 3234.0833
$ws.Range("I76").Value = 3258.842
$ws.Range("J76").Value = 3140
$ws.Range("K76").Value = 3258.842
$ws.Range("L76").Value = 3140
$ws.Range("M76").Value = -2943.842
$ws.Range("N76").Value = -3770
$ws.Range("H79").Value = 3234.0833
$ws.Range("I79").Value = 3258.842
$ws.Range("J79").Value = 3140
$ws.Range("K79").Value = 3258.842
$ws.Range("L79").Value = 3140
$ws.Range("M79").Value = -2166.842
$ws.Range("N79").Value = -5324
$ws.Range("H86").Value = 2727.6875
$ws.Range("I86").Value = 1606.8889
$ws.Range("J86").Value = 8780
$ws.Range("K86").Value = 1606.8889
$ws.Range("L86").Value = 8780
$ws.Range("M86").Value = -483.8888999999999
$ws.Range("N86").Value = -11026
$ws.Range("H89").Value = 2727.6875
$ws.Range("I89").Value = 1606.8889
$ws.Range("J89").Value = 8780
$ws.Range("K89").Value = 8034.4445
$ws.Range("L89").Value = 43900
$ws.Range("M89").Value = -2418.4445
$ws.Range("N89").Value = -55132
$ws.Range("H137").Value = 18183182
$ws.Range("I137").Value = 919.7045000000001
$ws.Range("K137").Value = 2759.1135
$ws.Range("M137").Value = -209.1135000000004

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1826.9487
$ws.Range("I61").Value = 1232.1613
$ws.Range("J61").Value = 4131.75
$ws.Range("K61").Value = 1232.1613
$ws.Range("L61").Value = 4131.75
$ws.Range("M61").Value = -1020.1613
$ws.Range("N61").Value = -4555.75
$ws.Range("H74").Value = 3996.675
$ws.Range("I74").Value = 711.2778
$ws.Range("J74").Value = 6684.727
$ws.Range("K74").Value = 711.2778
$ws.Range("L74").Value = 6684.727
$ws.Range("M74").Value = 162.7222
$ws.Range("N74").Value = -8432.726999999999
$ws.Range("H77").Value = 3996.675
$ws.Range("I77").Value = 711.2778
$ws.Range("J77").Value = 6684.727
$ws.Range("K77").Value = 3556.389
$ws.Range("L77").Value = 33423.635
$ws.Range("M77").Value = 811.6110000000003
$ws.Range("N77").Value = -42159.635
$ws.Range("H102").Value = 1386.6666
$ws.Range("I102").Value = 1386.6666
$ws.Range("K102").Value = 1386.6666
$ws.Range("M102").Value = 235.3334
$ws.Range("H132").Value = 1948.3243
$ws.Range("I132").Value = 1786.3334
$ws.Range("J132").Value = 2642.5715
$ws.Range("K132").Value = 5359.0002
$ws.Range("L132").Value = 7927.7145
$ws.Range("M132").Value = -2829.0002
$ws.Range("N132").Value = -12987.7145
$ws.Range("H136").Value = 1826.9487
$ws.Range("I136").Value = 1232.1613
$ws.Range("J136").Value = 4131.75
$ws.Range("K136").Value = 3696.4839
$ws.Range("L136").Value = 12395.25
$ws.Range("M136").Value = -1146.4839
$ws.Range("N136").Value = -17495.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 41008.77
$ws.Range("I134").Value = 52009.3
$ws.Range("J134").Value = 4340.3335
$ws.Range("K134").Value = 156027.9
$ws.Range("L134").Value = 13021.0005
$ws.Range("M134").Value = -153492.9
$ws.Range("N134").Value = -18091.0005

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1223.8
$ws.Range("I31").Value = 1147.76
$ws.Range("J31").Value = 1604
$ws.Range("K31").Value = 1147.76
$ws.Range("L31").Value = 1604
$ws.Range("M31").Value = -852.76
$ws.Range("N31").Value = -2194
$ws.Range("H34").Value = 1223.8
$ws.Range("I34").Value = 1147.76
$ws.Range("J34").Value = 1604
$ws.Range("K34").Value = 1147.76
$ws.Range("L34").Value = 1604
$ws.Range("M34").Value = -945.76
$ws.Range("N34").Value = -2008
$ws.Range("H58").Value = 1720.9231
$ws.Range("I58").Value = 1330
$ws.Range("K58").Value = 1330
$ws.Range("M58").Value = -1127
$ws.Range("H132").Value = 2335.1333
$ws.Range("I132").Value = 1820.875
$ws.Range("J132").Value = 2922.8572
$ws.Range("K132").Value = 5462.625
$ws.Range("L132").Value = 8768.571599999999
$ws.Range("M132").Value = -2932.625
$ws.Range("N132").Value = -13828.5716
$ws.Range("H134").Value = 2191.4443
$ws.Range("I134").Value = 2090.375
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 6271.125
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -3736.125
$ws.Range("N134").Value = -14070
$ws.Range("H136").Value = 1720.9231
$ws.Range("I136").Value = 1330
$ws.Range("K136").Value = 3990
$ws.Range("M136").Value = -1440

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1085.0834
$ws.Range("I69").Value = 752.625
$ws.Range("J69").Value = 1750
$ws.Range("K69").Value = 2257.875
$ws.Range("L69").Value = 5250
$ws.Range("M69").Value = -1446.875
$ws.Range("N69").Value = -6872
$ws.Range("H72").Value = 1085.0834
$ws.Range("I72").Value = 752.625
$ws.Range("J72").Value = 1750
$ws.Range("K72").Value = 6773.625
$ws.Range("L72").Value = 15750
$ws.Range("M72").Value = -2717.625
$ws.Range("N72").Value = -23862
$ws.Range("H95").Value = 3166.6667
$ws.Range("J95").Value = 3166.6667
$ws.Range("L95").Value = 9500.000100000001
$ws.Range("N95").Value = -13618.0001
$ws.Range("H100").Value = 2745
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 9000
$ws.Range("N100").Value = -10622
$ws.Range("H105").Value = 302668500
$ws.Range("J105").Value = 302668500
$ws.Range("L105").Value = 908005500
$ws.Range("N105").Value = -908010742
$ws.Range("H106").Value = 5500
$ws.Range("J106").Value = 5500
$ws.Range("L106").Value = 16500
$ws.Range("N106").Value = -18392
$ws.Range("H117").Value = 5024.25
$ws.Range("I117").Value = 495
$ws.Range("J117").Value = 5671.2856
$ws.Range("K117").Value = 1485
$ws.Range("L117").Value = 17013.8568
$ws.Range("M117").Value = 1957
$ws.Range("N117").Value = -23897.8568
$ws.Range("H120").Value = 6370
$ws.Range("I120").Value = 6777.5
$ws.Range("K120").Value = 20332.5
$ws.Range("M120").Value = -15494.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2364.4546
$ws.Range("I132").Value = 1001.8333
$ws.Range("J132").Value = 3999.6
$ws.Range("K132").Value = 3005.4999
$ws.Range("L132").Value = 11998.8
$ws.Range("M132").Value = -475.4998999999998
$ws.Range("N132").Value = -17058.8

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1470.1666
$ws.Range("I82").Value = 1251.909
$ws.Range("J82").Value = 1813.1428
$ws.Range("K82").Value = 1251.909
$ws.Range("L82").Value = 1813.1428
$ws.Range("M82").Value = -890.9090000000001
$ws.Range("N82").Value = -2535.1428
$ws.Range("H85").Value = 1470.1666
$ws.Range("I85").Value = 1251.909
$ws.Range("J85").Value = 1813.1428
$ws.Range("K85").Value = 1251.909
$ws.Range("L85").Value = 1813.1428
$ws.Range("M85").Value = -3.909000000000106
$ws.Range("N85").Value = -4309.1428
$ws.Range("H132").Value = 2697.0625
$ws.Range("I132").Value = 1647.7894
$ws.Range("J132").Value = 4230.615
$ws.Range("K132").Value = 4943.3682
$ws.Range("L132").Value = 12691.845
$ws.Range("M132").Value = -2413.3682
$ws.Range("N132").Value = -17751.845
$ws.Range("H136").Value = 2269
$ws.Range("I136").Value = 1427.9131
$ws.Range("J136").Value = 3558.6667
$ws.Range("K136").Value = 4283.7393
$ws.Range("L136").Value = 10676.0001
$ws.Range("M136").Value = -1733.7393
$ws.Range("N136").Value = -15776.0001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3515.2964
$ws.Range("I132").Value = 3745.5
$ws.Range("J132").Value = 2857.5715
$ws.Range("K132").Value = 11236.5
$ws.Range("L132").Value = 8572.7145
$ws.Range("M132").Value = -8706.5
$ws.Range("N132").Value = -13632.7145
$ws.Range("H136").Value = 12905.941
$ws.Range("I136").Value = 15457.214
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 46371.642
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -43821.642
$ws.Range("N136").Value = -8100
